$d = $word.ActiveDocument
$wordXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wordXmlNs + '><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($pkgXml)
}

# 1. Insert a blank paragraph right after the "Description: Start planning..." paragraph.
$descPlanning = $d.Paragraphs.Item(10)
$rngAfterPlanning = $descPlanning.Range
$rngAfterPlanning.Collapse(0)
$rngAfterPlanning.InsertParagraphAfter()

# 2. At the very end of the document (after the final "Description:" paragraph),
#    append 4 blank paragraphs that will become:
#      - a blank separator paragraph
#      - "9/27: Research / Defining Project Scope"
#      - "Start 9:00AM" (which also picks up the relocated _GoBack bookmark)
#      - a trailing blank paragraph
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$baseCount = $d.Paragraphs.Count

# Fill in the "9/27: Research / Defining Project Scope" paragraph, kept as two
# separate runs to mirror how the text was originally authored.
$pScope = $d.Paragraphs.Item($baseCount - 2)
Set-ParagraphXml $pScope '<w:r><w:t xml:space="preserve">9/27: Research </w:t></w:r><w:r><w:t>/ Defining Project Scope</w:t></w:r>'

# Fill in the "Start 9:00AM" paragraph and give it the _GoBack bookmark (moved from
# the "Start: 11:30 AM" paragraph further up in the document).
$pStart = $d.Paragraphs.Item($baseCount - 1)
Set-ParagraphXml $pStart '<w:r><w:t>Start 9:00AM</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

# 3. Remove the _GoBack bookmark from its original location on "Start: 11:30 AM".
$pOldBookmark = $d.Paragraphs.Item(7)
Set-ParagraphXml $pOldBookmark '<w:r><w:t>Start: 11:30 AM</w:t></w:r>'
